$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 0.001
$ws.Range("K4").Value = 684
$ws.Range("L4").Value = 0.001368
